$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a "Price" column cell (D) while forcing it to stay
# text (many of these look numeric, e.g. "1.000" or "334.00", and Excel would
# otherwise silently coerce them to numbers and drop the formatting).
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
}

# Row 2
Set-TextValue $ws.Range("D2") '30.803.84'
$ws.Range("E2").Value = '  +2.19%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.119.74'
$ws.Range("E3").Value = '  +10.37%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.000'
$ws.Range("E4").Value = '  -0.20%  '

# Row 5
Set-TextValue $ws.Range("D5") '334.00'
$ws.Range("E5").Value = '  +4.66%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.9988'
$ws.Range("E6").Value = '  -0.24%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.5237'

# Row 8
Set-TextValue $ws.Range("D8") '0.4411'
$ws.Range("E8").Value = '  +8.36%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.09064'
$ws.Range("E9").Value = '  +8.48%  '

# Row 10
Set-TextValue $ws.Range("D10") '46.77'
$ws.Range("E10").Value = '  +10.59%  '

# Row 11
Set-TextValue $ws.Range("D11") '1.187'
$ws.Range("E11").Value = '  +6.84%  '

# Row 12
Set-TextValue $ws.Range("D12") '25.38'
$ws.Range("E12").Value = '  +5.86%  '

# Row 13
Set-TextValue $ws.Range("D13") '2.116.52'
$ws.Range("E13").Value = '  +10.05%  '

# Row 14
Set-TextValue $ws.Range("D14") '6.769'
$ws.Range("E14").Value = '  +5.37%  '

# Row 15
Set-TextValue $ws.Range("D15") '7.864'
$ws.Range("E15").Value = '  +8.30%  '

# Row 16
Set-TextValue $ws.Range("D16") '98.24'
$ws.Range("E16").Value = '  +6.12%  '

# Row 17
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D17") '1.001'
$ws.Range("E17").Value = '  -0.57%  '

# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D18") '0.00001139'
$ws.Range("E18").Value = '  +3.98%  '

# Row 19
Set-TextValue $ws.Range("D19") '0.06633'
$ws.Range("E19").Value = '  +1.86%  '

# Row 20
Set-TextValue $ws.Range("D20") '19.22'
$ws.Range("E20").Value = '  +3.86%  '

# Row 21
Set-TextValue $ws.Range("D21") '6.406'
$ws.Range("E21").Value = '  +7.61%  '

# Row 22
Set-TextValue $ws.Range("D22") '0.9993'
$ws.Range("E22").Value = '  -0.22%  '

# Row 23
Set-TextValue $ws.Range("D23") '30.926.70'
$ws.Range("E23").Value = '  +2.57%  '

# Row 24
Set-TextValue $ws.Range("D24") '12.17'
$ws.Range("E24").Value = '  +7.09%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.367.31'
$ws.Range("E25").Value = '  +10.58%  '

# Row 26
Set-TextValue $ws.Range("D26") '2.257'
$ws.Range("E26").Value = '  +2.86%  '

# Row 27
Set-TextValue $ws.Range("D27") '22.98'
$ws.Range("E27").Value = '  +4.89%  '

# Row 28
Set-TextValue $ws.Range("D28") '2.582'

# Row 29
Set-TextValue $ws.Range("D29") '163.47'
$ws.Range("E29").Value = '  +0.51%  '

# Row 30
Set-TextValue $ws.Range("D30") '133.63'
$ws.Range("E30").Value = '  +3.62%  '

# Row 31
Set-TextValue $ws.Range("D31") '1.181'
$ws.Range("E31").Value = '  +4.54%  '

# Row 32
$ws.Range("E32").Value = '  +2.38%  '

# Row 33
Set-TextValue $ws.Range("D33") '6.253'
$ws.Range("E33").Value = '  +4.91%  '

# Row 34
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D34") '1.580'
$ws.Range("E34").Value = '  +31.48%  '

# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D35") '4.016'
$ws.Range("E35").Value = '  +5.80%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.02599'
$ws.Range("E36").Value = '  +5.81%  '

# Row 37
Set-TextValue $ws.Range("D37") '5.591'
$ws.Range("E37").Value = '  +4.97%  '

# Row 38
Set-TextValue $ws.Range("D38") '9.601'
$ws.Range("E38").Value = '  +11.77%  '

# Row 39
Set-TextValue $ws.Range("D39") '0.06764'
$ws.Range("E39").Value = '  +5.02%  '

# Row 40
Set-TextValue $ws.Range("D40") '12.72'
$ws.Range("E40").Value = '  +11.06%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.2269'
$ws.Range("E41").Value = '  +5.54%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.6846'

# Row 43
Set-TextValue $ws.Range("D43") '1.260'
$ws.Range("E43").Value = '  +4.19%  '

# Row 44
Set-TextValue $ws.Range("D44") '14.20'
$ws.Range("E44").Value = '  +5.87%  '

# Row 45
Set-TextValue $ws.Range("D45") '0.6439'
$ws.Range("E45").Value = '  +6.39%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.9986'
$ws.Range("E46").Value = '  +0.17%  '

# Row 47
Set-TextValue $ws.Range("D47") '2.262'
$ws.Range("E47").Value = '  +3.47%  '

# Row 48
Set-TextValue $ws.Range("D48") '3.680'
$ws.Range("E48").Value = '  +1.57%  '

# Row 49
Set-TextValue $ws.Range("D49") '1.286'
$ws.Range("E49").Value = '  +6.10%  '

# Row 50
Set-TextValue $ws.Range("D50") '83.18'
$ws.Range("E50").Value = '  +5.25%  '

# Row 51
$ws.Range("E51").Value = '  +3.62%  '
